$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '245.49'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '23.07'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.398'

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06050'

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8077'

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9310'

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07444'

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.03338'

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03068'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09375'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001599'

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.04827'

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.005256'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.004165'

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0009820'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.446'

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03981'

$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006433'
$ws.Range("E41").Value = '40KickTokenKICK'

$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1076'
$ws.Range("E42").Value = '41BKEXTokenBKK'

$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002901'
$ws.Range("E43").Value = '42CEJICEJI'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.006300'

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.9004'

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.002180'
